$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.814.93'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.636.34'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.18'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.629.93'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('E8').Value = '  +1.32%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('E10').Value = '  -5.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.80'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +16.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.610'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '48.58'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000285'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.221.01'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '675.03'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.69%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '8.97'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.634.08'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.840.59'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.29%  '
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('E21').Value = '  -4.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.51'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.88%  '
$ws.Range('E23').Value = '  +0.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '17.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '99.97'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.57%  '
$ws.Range('E26').Value = '  -2.70%  '
$ws.Range('E27').Value = '  -2.58%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E29').Value = '  -2.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.65'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.16'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.53'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.96%  '
$ws.Range('E34').Value = '  -6.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.98'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '575.99'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.35%  '
$ws.Range('E37').Value = '  -2.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.108'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '58.50'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0452'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.554.96'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.34%  '
$ws.Range('E43').Value = '  -1.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '34.37'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₃0734'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.69'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.70%  '
$ws.Range('E48').Value = '  +4.24%  '
$ws.Range('E49').Value = '  +1.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '137.79'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.91%  '
$ws.Range('E51').Value = '  -4.18%  '
